# Updates the worksheet date and the 25 division problems to the next
# day's randomly generated set, per commit "Update master to output
# generated at c8c62b6".

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title / date line
Replace-Text "2025-11-15 Saturday" "2025-11-16 Sunday"

# Row 1
Replace-Text "98÷7=" "87÷9="
Replace-Text "79÷9=" "76÷7="
Replace-Text "91÷8=" "18÷4="
Replace-Text "82÷4=" "96÷9="
Replace-Text "46÷2=" "56÷8="

# Row 2
Replace-Text "20÷9=" "14÷7="
Replace-Text "93÷3=" "69÷8="
Replace-Text "12÷9=" "98÷9="
Replace-Text "19÷6=" "25÷6="
Replace-Text "70÷3=" "39÷9="

# Row 3 - note "69÷7=" is both an old value (becomes "55÷6=") and the new
# value produced by another cell ("60÷5=" becomes "69÷7="). Replace the
# former before the latter so the freshly-written "69÷7=" is not re-matched.
Replace-Text "69÷7=" "55÷6="
Replace-Text "63÷3=" "27÷4="
Replace-Text "95÷4=" "45÷2="
Replace-Text "60÷5=" "69÷7="
Replace-Text "72÷7=" "21÷6="

# Row 4
Replace-Text "98÷8=" "22÷8="
Replace-Text "88÷2=" "35÷8="
Replace-Text "54÷9=" "35÷5="
Replace-Text "24÷7=" "21÷8="
Replace-Text "29÷4=" "88÷3="

# Row 5
Replace-Text "57÷9=" "85÷7="
Replace-Text "71÷9=" "41÷3="
Replace-Text "70÷9=" "65÷4="
Replace-Text "90÷6=" "67÷6="
Replace-Text "83÷9=" "55÷9="
